$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 952752.2
$ws.Range("I15").Value = 952752.2
$ws.Range("K15").Value = 2858256.6
$ws.Range("M15").Value = -2858087.6
$ws.Range("H41").Value = 899.375
$ws.Range("I41").Value = 1559.75
$ws.Range("K41").Value = 1559.75
$ws.Range("M41").Value = -1119.75
$ws.Range("H62").Value = 1698
$ws.Range("I62").Value = 1699.5
$ws.Range("J62").Value = 1695
$ws.Range("K62").Value = 1699.5
$ws.Range("L62").Value = 1695
$ws.Range("M62").Value = -1075.5
$ws.Range("N62").Value = -2943
$ws.Range("H65").Value = 1698
$ws.Range("I65").Value = 1699.5
$ws.Range("J65").Value = 1695
$ws.Range("K65").Value = 8497.5
$ws.Range("L65").Value = 8475
$ws.Range("M65").Value = -5377.5
$ws.Range("N65").Value = -14715
$ws.Range("H74").Value = 1504630.6
$ws.Range("I74").Value = 1504630.6
$ws.Range("K74").Value = 1504630.6
$ws.Range("M74").Value = -1503694.6
$ws.Range("H76").Value = 2721
$ws.Range("I76").Value = 2904.8
$ws.Range("K76").Value = 2904.8
$ws.Range("M76").Value = -2589.8
$ws.Range("H77").Value = 1504630.6
$ws.Range("I77").Value = 1504630.6
$ws.Range("K77").Value = 7523153
$ws.Range("M77").Value = -7518473
$ws.Range("H79").Value = 2721
$ws.Range("I79").Value = 2904.8
$ws.Range("K79").Value = 2904.8
$ws.Range("M79").Value = -1812.8
$ws.Range("H86").Value = 6823.231
$ws.Range("I86").Value = 7125.375
$ws.Range("J86").Value = 6339.8
$ws.Range("K86").Value = 7125.375
$ws.Range("L86").Value = 6339.8
$ws.Range("M86").Value = -6002.375
$ws.Range("N86").Value = -8585.799999999999
$ws.Range("H89").Value = 6823.231
$ws.Range("I89").Value = 7125.375
$ws.Range("J89").Value = 6339.8
$ws.Range("K89").Value = 35626.875
$ws.Range("L89").Value = 31699
$ws.Range("M89").Value = -30010.875
$ws.Range("N89").Value = -42931
$ws.Range("H95").Value = 39431.8
$ws.Range("J95").Value = 39431.8
$ws.Range("L95").Value = 39431.8
$ws.Range("N95").Value = -44923.8
$ws.Range("H115").Value = 287.5
$ws.Range("I115").Value = 287.5
$ws.Range("K115").Value = 862.5
$ws.Range("M115").Value = 704.5
$ws.Range("H116").Value = 7816.421
$ws.Range("I116").Value = 7816.421
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 7816.421
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -4374.421
$ws.Range("N116").ClearContents()
$ws.Range("H135").Value = 21739464
$ws.Range("I135").Value = 22727606
$ws.Range("J135").Value = 330
$ws.Range("K135").Value = 204548454
$ws.Range("L135").Value = 2970
$ws.Range("M135").Value = -204545919
$ws.Range("N135").Value = -8040

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 47621780
$ws.Range("I61").Value = 55557308
$ws.Range("J61").Value = 8596.333000000001
$ws.Range("K61").Value = 55557308
$ws.Range("L61").Value = 8596.333000000001
$ws.Range("M61").Value = -55557096
$ws.Range("N61").Value = -9020.333000000001
$ws.Range("H122").Value = 7771.8184
$ws.Range("I122").Value = 6686.25
$ws.Range("K122").Value = 20058.75
$ws.Range("M122").Value = -17608.75
$ws.Range("H132").Value = 6251979
$ws.Range("I132").Value = 7694456
$ws.Range("J132").Value = 1246.6666
$ws.Range("K132").Value = 23083368
$ws.Range("L132").Value = 3739.9998
$ws.Range("M132").Value = -23080838
$ws.Range("N132").Value = -8799.9998
$ws.Range("H136").Value = 47621780
$ws.Range("I136").Value = 55557308
$ws.Range("J136").Value = 8596.333000000001
$ws.Range("K136").Value = 166671924
$ws.Range("L136").Value = 25788.999
$ws.Range("M136").Value = -166669374
$ws.Range("N136").Value = -30888.999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 13638
$ws.Range("I75").Value = 11398.8
$ws.Range("J75").Value = 19236
$ws.Range("K75").Value = 11398.8
$ws.Range("L75").Value = 19236
$ws.Range("M75").Value = -10462.8
$ws.Range("N75").Value = -21108
$ws.Range("H78").Value = 13638
$ws.Range("I78").Value = 11398.8
$ws.Range("J78").Value = 19236
$ws.Range("K78").Value = 34196.39999999999
$ws.Range("L78").Value = 57708
$ws.Range("M78").Value = -29516.39999999999
$ws.Range("N78").Value = -67068
$ws.Range("H94").Value = 3444.0908
$ws.Range("J94").Value = 3000
$ws.Range("L94").Value = 3000
$ws.Range("N94").Value = -3902
$ws.Range("H134").Value = 15459656
$ws.Range("I134").Value = 17005254
$ws.Range("K134").Value = 51015762
$ws.Range("M134").Value = -51013227

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3180.44
$ws.Range("I31").Value = 3271.4707
$ws.Range("K31").Value = 3271.4707
$ws.Range("M31").Value = -2976.4707
$ws.Range("H34").Value = 3180.44
$ws.Range("I34").Value = 3271.4707
$ws.Range("K34").Value = 3271.4707
$ws.Range("M34").Value = -3069.4707
$ws.Range("H54").Value = 32749.25
$ws.Range("J54").Value = 32749.25
$ws.Range("L54").Value = 32749.25
$ws.Range("N54").Value = -34065.25
$ws.Range("H58").Value = 31258548
$ws.Range("I58").Value = 100021200
$ws.Range("K58").Value = 100021200
$ws.Range("M58").Value = -100020997
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H94").Value = 2483.375
$ws.Range("I94").Value = 2654.4
$ws.Range("K94").Value = 2654.4
$ws.Range("M94").Value = -2203.4
$ws.Range("H132").Value = 40002550
$ws.Range("I132").Value = 47621444
$ws.Range("K132").Value = 142864332
$ws.Range("M132").Value = -142861802
$ws.Range("H135").Value = 120000
$ws.Range("J135").Value = 120000
$ws.Range("L135").Value = 120000
$ws.Range("N135").Value = -130140
$ws.Range("H136").Value = 31258548
$ws.Range("I136").Value = 100021200
$ws.Range("K136").Value = 300063600
$ws.Range("M136").Value = -300061050

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 784.75
$ws.Range("I98").Value = 784.75
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2354.25
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -856.25
$ws.Range("N98").ClearContents()
$ws.Range("H121").Value = 126307.6
$ws.Range("I121").Value = 338809.66
$ws.Range("J121").Value = 35235.285
$ws.Range("K121").Value = 1016428.98
$ws.Range("L121").Value = 105705.855
$ws.Range("M121").Value = -1015118.98
$ws.Range("N121").Value = -108325.855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3199.875
$ws.Range("I22").Value = 3199.875
$ws.Range("K22").Value = 3199.875
$ws.Range("M22").Value = -2904.875
$ws.Range("H27").Value = 3199.875
$ws.Range("I27").Value = 3199.875
$ws.Range("K27").Value = 3199.875
$ws.Range("M27").Value = -3092.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 8003.3335
$ws.Range("I29").Value = 7005
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 7005
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = -6715
$ws.Range("N29").Value = -10580
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4708
$ws.Range("H94").Value = 19999
$ws.Range("J94").Value = 19999
$ws.Range("L94").Value = 19999
$ws.Range("N94").Value = -21801
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 17243152
$ws.Range("I136").Value = 21740894
$ws.Range("K136").Value = 65222682
$ws.Range("M136").Value = -65220132
$ws.Range("H140").Value = 74992
$ws.Range("J140").Value = 74992
$ws.Range("L140").Value = 74992
$ws.Range("N140").Value = -85352
$ws.Range("H141").Value = 200048800
$ws.Range("J141").Value = 200048800
$ws.Range("L141").Value = 200048800
$ws.Range("N141").Value = -200059160
